$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Target OOXML column widths (29.9777047293527 and 13.7470528738839) are not
# exactly reachable through the ColumnWidth COM property, which snaps to a
# 1/6-character pixel grid. 29.17 / 12.83 are the calibrated inputs that land
# on the closest achievable grid points (30 and 13.666666666666666).
$wideColWidth = 29.17
$narrowColWidth = 12.83

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-22 12:51:34"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(5).ColumnWidth = $wideColWidth
$wsZhCn.Columns.Item(6).ColumnWidth = $wideColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-22 12:51:41"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $narrowColWidth
